# Update "想去人数" (interest/attendance counts) figures in the "展览" and
# "全部类型" sheets to reflect newly scraped data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3303
$ws1.Range("F4").Value = 125
$ws1.Range("F5").Value = 6923
$ws1.Range("F6").Value = 2239
$ws1.Range("F9").Value = 26
$ws1.Range("F13").Value = 159
$ws1.Range("F14").Value = 249

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3303
$ws4.Range("F5").Value = 125
$ws4.Range("F6").Value = 6923
$ws4.Range("F7").Value = 2239
$ws4.Range("F10").Value = 26
$ws4.Range("F14").Value = 159
$ws4.Range("F15").Value = 249
